$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (Voltage), shifting Voltage..Library Path right by one.
$ws.Columns("F").EntireColumn.Insert()

# Set the new column's width to match the authored width (~9.5703125 characters).
$ws.Columns("F").ColumnWidth = 8.666666666666666

# Header for the newly inserted "Dielectric" column.
$ws.Range("F1").Value = "Dielectric"

# Dielectric values per row, matching each part's description (X7R vs X5R).
$ws.Range("F2:F7").Value = "X7R"
$ws.Range("F8:F11").Value = "X5R"
